# Update gh-pages to output generated at 456a3b4
# Applies updated "F" column (review/comment count) values across the
# four worksheets, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param($SheetName, $RowToValue)
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowToValue.Keys) {
        $ws.Range("F$row").Value = $RowToValue[$row]
    }
}

# 展览 sheet
Set-FValues "展览" @{
    4  = 64
    5  = 75
    6  = 848
    7  = 425
    8  = 4719
    9  = 4719
    11 = 123
    12 = 161
    15 = 120
    16 = 7539
    21 = 525
    22 = 1379
    24 = 6286
    25 = 2236
    26 = 20
    27 = 2090
    28 = 6182
    29 = 143
    33 = 448
    34 = 6447
    36 = 207
    39 = 21
    40 = 16
    41 = 2456
    43 = 61
    44 = 1118
    46 = 444
    47 = 2150
    48 = 47
    49 = 1078
}

# 演出 sheet
Set-FValues "演出" @{
    3  = 233
    6  = 125
    14 = 24
}

# 本地生活 sheet
Set-FValues "本地生活" @{
    2 = 1448
}

# 全部类型 sheet
Set-FValues "全部类型" @{
    3  = 1448
    4  = 64
    5  = 233
    6  = 75
    8  = 425
    9  = 4719
    10 = 4719
    12 = 123
    13 = 161
    16 = 120
    17 = 7539
    20 = 525
    21 = 1379
    22 = 125
    23 = 6286
    24 = 2239
    26 = 2090
    29 = 6182
    30 = 143
    35 = 448
    36 = 6447
    38 = 207
    40 = 21
    42 = 2456
    44 = 1118
    46 = 444
    48 = 2150
    49 = 47
    50 = 24
}
